$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")

# Update the admin test user's username to "metCs" (password "MetCs673" stays the same).
$ws.Range("C4").Value = "metCs"
$ws.Range("D4").Value = "MetCs673"

# Update the active selection to reflect where the user last clicked.
$ws.Range("C4").Select()
